$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.029816075909231
$ws.Range("D2").Value = 1.034624077967305
$ws.Range("E2").Value = 1.029600393852309
$ws.Range("F2").Value = 1.040027720571122
$ws.Range("I2").Value = 1.035887974341571
$ws.Range("J2").Value = 1.034960759601254
$ws.Range("K2").Value = 1.037423018030271
$ws.Range("L2").Value = 1.032413834043669
$ws.Range("M2").Value = 1.042811233608109
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030684986950767
$ws.Range("D3").Value = 1.035282926524373
$ws.Range("E3").Value = 1.030336523255732
$ws.Range("F3").Value = 1.041150918177463
$ws.Range("I3").Value = 1.03610282454456
$ws.Range("J3").Value = 1.035471274015288
$ws.Range("K3").Value = 1.037891493489876
$ws.Range("L3").Value = 1.032958336787458
$ws.Range("M3").Value = 1.043743946924692
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031247696985664
$ws.Range("D4").Value = 1.035709605057179
$ws.Range("E4").Value = 1.030813627786324
$ws.Range("F4").Value = 1.041878544495064
$ws.Range("I4").Value = 1.036240899560263
$ws.Range("J4").Value = 1.035801437902029
$ws.Range("K4").Value = 1.038194304356032
$ws.Range("L4").Value = 1.033310779432933
$ws.Range("M4").Value = 1.044347733295685
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031484370762737
$ws.Range("D5").Value = 1.035889065324612
$ws.Range("E5").Value = 1.031014387682581
$ws.Range("F5").Value = 1.042184639704969
$ws.Range("I5").Value = 1.036298718935356
$ws.Range("J5").Value = 1.03594019639548
$ws.Range("K5").Value = 1.03832152745851
$ws.Range("L5").Value = 1.033458972274852
$ws.Range("M5").Value = 1.044601626359127
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031524115775026
$ws.Range("D6").Value = 1.035919202412895
$ws.Range("E6").Value = 1.03104810696682
$ws.Range("F6").Value = 1.042236046214531
$ws.Range("I6").Value = 1.036308413714482
$ws.Range("J6").Value = 1.035963492042305
$ws.Range("K6").Value = 1.038342884139875
$ws.Range("L6").Value = 1.033483855995369
$ws.Range("M6").Value = 1.044644259697178
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031250858999645
$ws.Range("D7").Value = 1.035712002683401
$ws.Range("E7").Value = 1.030816309624225
$ws.Range("F7").Value = 1.041882633764712
$ws.Range("I7").Value = 1.036241673039634
$ws.Range("J7").Value = 1.035803292167643
$ws.Range("K7").Value = 1.038196004627393
$ws.Range("L7").Value = 1.033312759492325
$ws.Range("M7").Value = 1.044351125587214
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.030109631151911
$ws.Range("D8").Value = 1.034846663625759
$ws.Range("E8").Value = 1.029849009922039
$ws.Range("F8").Value = 1.0404071357078
$ws.Range("I8").Value = 1.035960779817802
$ws.Range("J8").Value = 1.035133325590935
$ws.Range("K8").Value = 1.037581408007334
$ws.Range("L8").Value = 1.032597827251174
$ws.Range("M8").Value = 1.04312639455098
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02810226765581
$ws.Range("D9").Value = 1.033324642390416
$ws.Range("E9").Value = 1.028150540373798
$ws.Range("F9").Value = 1.037813595627136
$ws.Range("I9").Value = 1.035458577953235
$ws.Range("J9").Value = 1.03395147614575
$ws.Range("K9").Value = 1.036495972299002
$ws.Range("L9").Value = 1.031338938527373
$ws.Range("M9").Value = 1.040970268550675
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026766533743689
$ws.Range("D10").Value = 1.032311940949756
$ws.Range("E10").Value = 1.027022373581495
$ws.Range("F10").Value = 1.036088954114567
$ws.Range("I10").Value = 1.035118946433982
$ws.Range("J10").Value = 1.033162771333962
$ws.Range("K10").Value = 1.03577076848386
$ws.Range("L10").Value = 1.030500357829066
$ws.Range("M10").Value = 1.039534234896497
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026188755791212
$ws.Range("D11").Value = 1.031873917470973
$ws.Range("E11").Value = 1.026534866776013
$ws.Range("F11").Value = 1.035343211312661
$ws.Range("I11").Value = 1.034970743242353
$ws.Range("J11").Value = 1.032821074380962
$ws.Range("K11").Value = 1.035456384802273
$ws.Range("L11").Value = 1.030137417513461
$ws.Range("M11").Value = 1.038912751401931
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02597423516669
$ws.Range("D12").Value = 1.031711290038052
$ws.Range("E12").Value = 1.026353936289237
$ws.Range("F12").Value = 1.035066365628131
$ws.Range("I12").Value = 1.034915523142878
$ws.Range("J12").Value = 1.032694126426367
$ws.Range("K12").Value = 1.035339554742589
$ws.Range("L12").Value = 1.030002632028814
$ws.Range("M12").Value = 1.038681954593115
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026020246394874
$ws.Range("D13").Value = 1.031746170803227
$ws.Range("E13").Value = 1.026392739624841
$ws.Range("F13").Value = 1.035125742883551
$ws.Range("I13").Value = 1.034927375764079
$ws.Range("J13").Value = 1.032721358377353
$ws.Range("K13").Value = 1.035364617613292
$ws.Range("L13").Value = 1.030031542747041
$ws.Range("M13").Value = 1.038731459029403
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026171021569496
$ws.Range("D14").Value = 1.031860473118963
$ws.Range("E14").Value = 1.026519907911588
$ws.Range("F14").Value = 1.035320323969264
$ws.Range("I14").Value = 1.034966182213724
$ws.Range("J14").Value = 1.032810581358449
$ws.Range("K14").Value = 1.03544672868829
$ws.Range("L14").Value = 1.030126275552229
$ws.Range("M14").Value = 1.038893672646398
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026263931333715
$ws.Range("D15").Value = 1.03193090839546
$ws.Range("E15").Value = 1.026598280576078
$ws.Range("F15").Value = 1.035440232534148
$ws.Range("I15").Value = 1.034990069526143
$ws.Range("J15").Value = 1.032865551099377
$ws.Range("K15").Value = 1.035497312910129
$ws.Range("L15").Value = 1.03018464713118
$ws.Range("M15").Value = 1.038993624427068
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026804891779017
$ws.Range("D16").Value = 1.032341021422767
$ws.Range("E16").Value = 1.02705474893168
$ws.Range("F16").Value = 1.036138468518781
$ws.Range("I16").Value = 1.035128758206313
$ws.Range("J16").Value = 1.033185444868836
$ws.Range("K16").Value = 1.035791625472624
$ws.Range("L16").Value = 1.030524448669343
$ws.Range("M16").Value = 1.03957548767154
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027144384392313
$ws.Range("D17").Value = 1.032598404797286
$ws.Range("E17").Value = 1.027341347404526
$ws.Range("F17").Value = 1.036576731778477
$ws.Range("I17").Value = 1.035215448967314
$ws.Range("J17").Value = 1.033386057561738
$ws.Range("K17").Value = 1.035976142759833
$ws.Range("L17").Value = 1.030737643641957
$ws.Range("M17").Value = 1.039940563406036
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027342462856222
$ws.Range("D18").Value = 1.032748578664415
$ws.Range("E18").Value = 1.027508611528277
$ws.Range("F18").Value = 1.036832463482204
$ws.Range("I18").Value = 1.035265904145086
$ws.Range("J18").Value = 1.03350305387373
$ws.Range("K18").Value = 1.036083733227159
$ws.Range("L18").Value = 1.030862013057122
$ws.Range("M18").Value = 1.040153537574125
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027410012324949
$ws.Range("D19").Value = 1.032799791915659
$ws.Range("E19").Value = 1.027565660524345
$ws.Range("F19").Value = 1.036919678361293
$ws.Range("I19").Value = 1.035283089370402
$ws.Range("J19").Value = 1.033542943565755
$ws.Range("K19").Value = 1.036120412771539
$ws.Range("L19").Value = 1.030904422561995
$ws.Range("M19").Value = 1.040226161619033
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02710795402315
$ws.Range("D20").Value = 1.032570785177029
$ws.Range("E20").Value = 1.027310588151677
$ws.Range("F20").Value = 1.036529699915532
$ws.Range("I20").Value = 1.035206159251109
$ws.Range("J20").Value = 1.03336453554008
$ws.Range("K20").Value = 1.035956349449992
$ws.Range("L20").Value = 1.030714768135552
$ws.Range("M20").Value = 1.039901390955656
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026126619494327
$ws.Range("D21").Value = 1.03182681188188
$ws.Range("E21").Value = 1.026482455833691
$ws.Range("F21").Value = 1.035263020368811
$ws.Range("I21").Value = 1.034954759393296
$ws.Range("J21").Value = 1.032784308138133
$ws.Range("K21").Value = 1.035422550505091
$ws.Range("L21").Value = 1.030098378356388
$ws.Range("M21").Value = 1.03884590339982
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025510147409958
$ws.Range("D22").Value = 1.031359475136302
$ws.Range("E22").Value = 1.025962651977069
$ws.Range("F22").Value = 1.034467513725503
$ws.Range("I22").Value = 1.034795706179134
$ws.Range("J22").Value = 1.03241934332975
$ws.Range("K22").Value = 1.035086618193413
$ws.Range("L22").Value = 1.029710984693766
$ws.Range("M22").Value = 1.038182565035536
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025836900127369
$ws.Range("D23").Value = 1.031607178063857
$ws.Range("E23").Value = 1.026238126326933
$ws.Range("F23").Value = 1.034889140925182
$ws.Range("I23").Value = 1.034880116790506
$ws.Range("J23").Value = 1.032612832252006
$ws.Range("K23").Value = 1.035264731491643
$ws.Range("L23").Value = 1.029916334349735
$ws.Range("M23").Value = 1.038534185742593
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027124415159226
$ws.Range("D24").Value = 1.032583265151041
$ws.Range("E24").Value = 1.027324486637313
$ws.Range("F24").Value = 1.0365509512804
$ws.Range("I24").Value = 1.035210357214156
$ws.Range("J24").Value = 1.033374260469977
$ws.Range("K24").Value = 1.035965293304025
$ws.Range("L24").Value = 1.03072510454169
$ws.Range("M24").Value = 1.039919091201939
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028620782295031
$ws.Range("D25").Value = 1.033717778741112
$ws.Range("E25").Value = 1.028588911091332
$ws.Range("F25").Value = 1.038483317284896
$ws.Range("I25").Value = 1.035589262906715
$ws.Range("J25").Value = 1.03425715866987
$ws.Range("K25").Value = 1.036776866265035
$ws.Range("L25").Value = 1.03166427657412
$ws.Range("M25").Value = 1.041527437474784
